$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LoginTestData")

# --- Update data rows 2-7 with new emails and a literal numeric password value ---
$emails = @("jakay11@gmail.com", "jakay12@gmail.com", "jakay13@gmail.com", "jakay14@gmail.com", "jakay15@gmail.com", "jakay16@gmail.com")
$modes  = @("y", "y", "n", "n", "n", "n")

for ($i = 0; $i -lt 6; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $emails[$i]
    $ws.Cells.Item($r, 2).Value = 12345678
    $ws.Cells.Item($r, 3).Value = $modes[$i]
}

# Row 8 keeps its original values (jakay34@gmail.com / 12345678 / y) - untouched.

# --- Hyperlinks: the underlying mailto targets got shuffled between rows (A5 and A8
# effectively trade hyperlink targets versus the untouched layout), so rebuild the
# whole collection in the exact order/target combination seen in the saved file. ---
$ws.Cells.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:automation@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A8"), "mailto:automation@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A6"), "mailto:automation@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A7"), "mailto:automation@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:automation@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:automation@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A5"), "mailto:jakay34@gmail.com")

# Adding a hyperlink re-stamps the cell with a fresh "Hyperlink" style entry;
# restore the original shared Hyperlink style index so column A keeps the same
# cellXfs it had before (bold/underline hyperlink look, no stray style bloat).
$ws.Range("A2").Style = "Hyperlink"
$ws.Range("A3").Style = "Hyperlink"
$ws.Range("A4").Style = "Hyperlink"
$ws.Range("A5").Style = "Hyperlink"
$ws.Range("A6").Style = "Hyperlink"
$ws.Range("A7").Style = "Hyperlink"
$ws.Range("A8").Style = "Hyperlink"

# --- Window width (best effort - engine may not persist this) ---
$excel.ActiveWindow.Width = 19815

# --- Update the active selection to D8 ---
$ws.Range("D8").Select()
